$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.64723265171051
$ws.Range("B1").Value = 3.612305402755737
$ws.Range("C1").Value = 3.222667455673218
$ws.Range("D1").Value = 2.60638952255249
$ws.Range("E1").Value = 1.662910103797913
